$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: refresh status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- zh-cn sheet: status + target/handback file + handback datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("F2").Value = "ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.md"
$wsZh.Range("G2").Value = "ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.00833d82b779c06f28bfaf2fe57da27d8779351e.zh-cn.xlf"
$wsZh.Range("F3").Value = "fecaeb59-b73e-4eee-898d-bd600c3b8478.md"
$wsZh.Range("G3").Value = "fecaeb59-b73e-4eee-898d-bd600c3b8478.4dfcd94b7a02e0f79e16d467dc88cdd6987d7526.zh-cn.xlf"

$wsZh.Range("H2").Value = "2016-03-17 20:11:28"
$wsZh.Range("H3").Value = "2016-03-17 20:11:28"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/b09fafa0151cda50dfc3d4624967cd4904691662/e2e/ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.md", [Type]::Missing, [Type]::Missing, "ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b842921d45d012d2a27a83056e07123eda6c71a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.00833d82b779c06f28bfaf2fe57da27d8779351e.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.00833d82b779c06f28bfaf2fe57da27d8779351e.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/b09fafa0151cda50dfc3d4624967cd4904691662/e2e/fecaeb59-b73e-4eee-898d-bd600c3b8478.md", [Type]::Missing, [Type]::Missing, "fecaeb59-b73e-4eee-898d-bd600c3b8478.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b842921d45d012d2a27a83056e07123eda6c71a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/fecaeb59-b73e-4eee-898d-bd600c3b8478.4dfcd94b7a02e0f79e16d467dc88cdd6987d7526.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "fecaeb59-b73e-4eee-898d-bd600c3b8478.4dfcd94b7a02e0f79e16d467dc88cdd6987d7526.zh-cn.xlf")

# --- de-de sheet: status + target/handback file + handback datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("F2").Value = "ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.md"
$wsDe.Range("G2").Value = "ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.00833d82b779c06f28bfaf2fe57da27d8779351e.de-de.xlf"
$wsDe.Range("F3").Value = "fecaeb59-b73e-4eee-898d-bd600c3b8478.md"
$wsDe.Range("G3").Value = "fecaeb59-b73e-4eee-898d-bd600c3b8478.4dfcd94b7a02e0f79e16d467dc88cdd6987d7526.de-de.xlf"

$wsDe.Range("H2").Value = "2016-03-17 20:11:34"
$wsDe.Range("H3").Value = "2016-03-17 20:11:34"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/b09fafa0151cda50dfc3d4624967cd4904691662/e2e/ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.md", [Type]::Missing, [Type]::Missing, "ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a330be5b780384daa3994fd18e55487859bb345/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.00833d82b779c06f28bfaf2fe57da27d8779351e.de-de.xlf", [Type]::Missing, [Type]::Missing, "ac728b7b-b3eb-4f35-bf0b-6e1c8ae53351.00833d82b779c06f28bfaf2fe57da27d8779351e.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/b09fafa0151cda50dfc3d4624967cd4904691662/e2e/fecaeb59-b73e-4eee-898d-bd600c3b8478.md", [Type]::Missing, [Type]::Missing, "fecaeb59-b73e-4eee-898d-bd600c3b8478.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a330be5b780384daa3994fd18e55487859bb345/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/fecaeb59-b73e-4eee-898d-bd600c3b8478.4dfcd94b7a02e0f79e16d467dc88cdd6987d7526.de-de.xlf", [Type]::Missing, [Type]::Missing, "fecaeb59-b73e-4eee-898d-bd600c3b8478.4dfcd94b7a02e0f79e16d467dc88cdd6987d7526.de-de.xlf")
